# Update the Assignment List worksheet:
#  - remove the old "Essay" assignment row (CMST 1234 Essay Fri May 10 2019)
#  - shift the remaining rows up
#  - append two new assignment rows at the bottom of the table
#    (Reflection Essay due 5/12, and a TEST "testing 123" due 6/12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row that holds the old "Essay" assignment (row 2);
# this shifts rows 3-5 up to become rows 2-4.
$ws.Rows.Item(2).Delete()

# New row 5: CMST / 1010 / Reflection Essay / Sun May 12 2019
$ws.Cells.Item(5, 1).Value = "CMST"
$ws.Cells.Item(5, 2).Value = 1010.0
$ws.Cells.Item(5, 3).Value = "Reflection Essay"
$ws.Cells.Item(5, 4).Value = "Sun May 12 00:00:00 CDT 2019"

# New row 6: TEST / 4523 / testing 123 / Wed Jun 12 2019
$ws.Cells.Item(6, 1).Value = "TEST"
$ws.Cells.Item(6, 2).Value = 4523.0
$ws.Cells.Item(6, 3).Value = "testing 123"
$ws.Cells.Item(6, 4).Value = "Wed Jun 12 00:00:00 CDT 2019"
